$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 247.61539
$ws.Range("I39").Value = 164.44444
$ws.Range("K39").Value = 493.33332
$ws.Range("M39").Value = -197.33332

$ws.Range("H107").Value = 561.4
$ws.Range("I107").Value = 612.7778
$ws.Range("K107").Value = 612.7778
$ws.Range("M107").Value = 1307.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2750
$ws.Range("I26").Value = 1000
$ws.Range("K26").Value = 1000
$ws.Range("M26").Value = -670

$ws.Range("H32").Value = 6112.067
$ws.Range("I32").Value = 4370.4614
$ws.Range("J32").Value = 17432.5
$ws.Range("K32").Value = 4370.4614
$ws.Range("L32").Value = 17432.5
$ws.Range("M32").Value = -4083.4614
$ws.Range("N32").Value = -18006.5

$ws.Range("H74").Value = 925.05884
$ws.Range("I74").Value = 925.05884
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 925.05884
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -51.05884000000003
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 925.05884
$ws.Range("I77").Value = 925.05884
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4625.2942
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -257.2942000000003
$ws.Range("N77").ClearContents()

$ws.Range("H88").Value = 1043.7142
$ws.Range("I88").Value = 433.33334
$ws.Range("J88").Value = 1501.5
$ws.Range("K88").Value = 433.33334
$ws.Range("L88").Value = 1501.5
$ws.Range("M88").Value = -27.33334000000002
$ws.Range("N88").Value = -2313.5

$ws.Range("H91").Value = 1043.7142
$ws.Range("I91").Value = 433.33334
$ws.Range("J91").Value = 1501.5
$ws.Range("K91").Value = 433.33334
$ws.Range("L91").Value = 1501.5
$ws.Range("M91").Value = 970.66666
$ws.Range("N91").Value = -4309.5

$ws.Range("H132").Value = 200
$ws.Range("I132").Value = 200
$ws.Range("K132").Value = 600
$ws.Range("M132").Value = 1930

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1916.6
$ws.Range("I105").Value = 1448.4166
$ws.Range("K105").Value = 1448.4166
$ws.Range("M105").Value = 298.5834

$ws.Range("H107").Value = 3917
$ws.Range("I107").Value = 4099.6665
$ws.Range("J107").Value = 3369
$ws.Range("K107").Value = 4099.6665
$ws.Range("L107").Value = 3369
$ws.Range("M107").Value = -2179.6665
$ws.Range("N107").Value = -7209

$ws.Range("H134").Value = 1883.0667
$ws.Range("I134").Value = 1432.125
$ws.Range("J134").Value = 3686.8333
$ws.Range("K134").Value = 4296.375
$ws.Range("L134").Value = 11060.4999
$ws.Range("M134").Value = -1761.375
$ws.Range("N134").Value = -16130.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5311.5
$ws.Range("I31").Value = 4967.8
$ws.Range("K31").Value = 4967.8
$ws.Range("M31").Value = -4672.8

$ws.Range("H34").Value = 5311.5
$ws.Range("I34").Value = 4967.8
$ws.Range("K34").Value = 4967.8
$ws.Range("M34").Value = -4765.8

$ws.Range("H58").Value = 1897.8368
$ws.Range("J58").Value = 4800
$ws.Range("L58").Value = 4800
$ws.Range("N58").Value = -5206

$ws.Range("H99").Value = 4113.7
$ws.Range("I99").Value = 3338.5
$ws.Range("J99").Value = 4630.5
$ws.Range("K99").Value = 3338.5
$ws.Range("L99").Value = 4630.5
$ws.Range("M99").Value = -1840.5
$ws.Range("N99").Value = -7626.5

$ws.Range("H105").Value = 1516.1666
$ws.Range("I105").Value = 1201.3334
$ws.Range("J105").Value = 1831
$ws.Range("K105").Value = 1201.3334
$ws.Range("L105").Value = 1831
$ws.Range("M105").Value = 545.6666
$ws.Range("N105").Value = -5325

$ws.Range("H122").Value = 3825.1
$ws.Range("I122").Value = 3825.1
$ws.Range("K122").Value = 11475.3
$ws.Range("M122").Value = -9025.299999999999

$ws.Range("H126").Value = 4113.7
$ws.Range("I126").Value = 3338.5
$ws.Range("J126").Value = 4630.5
$ws.Range("K126").Value = 10015.5
$ws.Range("L126").Value = 13891.5
$ws.Range("M126").Value = -7545.5
$ws.Range("N126").Value = -18831.5

$ws.Range("H136").Value = 1897.8368
$ws.Range("J136").Value = 4800
$ws.Range("L136").Value = 14400
$ws.Range("N136").Value = -19500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 55454780
$ws.Range("J4").Value = 617.61536
$ws.Range("K4").Value = 166364340
$ws.Range("L4").Value = 1852.84608
$ws.Range("M4").Value = -166364228
$ws.Range("N4").Value = -2076.84608

$ws.Range("H107").Value = 1643
$ws.Range("J107").Value = 1831.75
$ws.Range("L107").Value = 5495.25
$ws.Range("N107").Value = -9335.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960

$ws.Range("H129").Value = 50000
$ws.Range("J129").Value = 50000
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000

$ws.Range("H132").Value = 2350.9167
$ws.Range("I132").Value = 2199.7222
$ws.Range("K132").Value = 6599.1666
$ws.Range("M132").Value = -4069.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2731.25
$ws.Range("I7").Value = 2621.4285
$ws.Range("K7").Value = 2621.4285
$ws.Range("M7").Value = -2509.4285

$ws.Range("H22").Value = 3789.818
$ws.Range("I22").Value = 2969.5715
$ws.Range("J22").Value = 5225.25
$ws.Range("K22").Value = 2969.5715
$ws.Range("L22").Value = 5225.25
$ws.Range("M22").Value = -2674.5715
$ws.Range("N22").Value = -5815.25

$ws.Range("H27").Value = 3789.818
$ws.Range("I27").Value = 2969.5715
$ws.Range("J27").Value = 5225.25
$ws.Range("K27").Value = 2969.5715
$ws.Range("L27").Value = 5225.25
$ws.Range("M27").Value = -2862.5715
$ws.Range("N27").Value = -5439.25

$ws.Range("H34").Value = 33500
$ws.Range("I34").Value = 33000
$ws.Range("K34").Value = 33000
$ws.Range("M34").Value = -32828

$ws.Range("H43").Value = 4067603.5
$ws.Range("I43").Value = 8000
$ws.Range("J43").Value = 5082504
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 5082504
$ws.Range("M43").Value = -7807
$ws.Range("N43").Value = -5082890

$ws.Range("H93").Value = 1790.5454
$ws.Range("I93").Value = 1525.75
$ws.Range("K93").Value = 1525.75
$ws.Range("M93").Value = -277.75

$ws.Range("H126").Value = 2731.25
$ws.Range("I126").Value = 2621.4285
$ws.Range("K126").Value = 7864.2855
$ws.Range("M126").Value = -5394.2855

$ws.Range("H132").Value = 4471.467
$ws.Range("I132").Value = 3091.6667
$ws.Range("J132").Value = 9990.666999999999
$ws.Range("K132").Value = 9275.000100000001
$ws.Range("L132").Value = 29972.001
$ws.Range("M132").Value = -6745.000100000001
$ws.Range("N132").Value = -35032.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 42499.8
$ws.Range("J48").Value = 42499.8
$ws.Range("L48").Value = 42499.8
$ws.Range("N48").Value = -43637.8

$ws.Range("H54").Value = 35724.145
$ws.Range("I54").Value = 33345
$ws.Range("J54").Value = 49999
$ws.Range("K54").Value = 33345
$ws.Range("L54").Value = 49999
$ws.Range("M54").Value = -32825
$ws.Range("N54").Value = -51039

$ws.Range("H117").Value = 120409
$ws.Range("J117").Value = 120409
$ws.Range("L117").Value = 120409
$ws.Range("N117").Value = -129587

$ws.Range("H132").Value = 2931.05
$ws.Range("I132").Value = 2573.3928
$ws.Range("J132").Value = 3765.5833
$ws.Range("K132").Value = 7720.178400000001
$ws.Range("L132").Value = 11296.7499
$ws.Range("M132").Value = -5190.178400000001
$ws.Range("N132").Value = -16356.7499

$ws.Range("H136").Value = 1712.625
$ws.Range("I136").Value = 1225.2354
$ws.Range("J136").Value = 2896.2856
$ws.Range("K136").Value = 3675.7062
$ws.Range("L136").Value = 8688.856800000001
$ws.Range("M136").Value = -1125.7062
$ws.Range("N136").Value = -13788.8568

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
